# "Balancing testing in Code"
#
# Sheet1: the "Balancing Request" detail column (B) for the BMC balancing-
#   request IDs (rows 66-83) is reshuffled - most of the old column-B values
#   slide over into a new column H (and the label "And Balancing state" is
#   added next to the CMC8 temperature row), while column B on a few rows
#   is rewritten / cleared as the CMC numbering is corrected.
#
# Sheet2: a third bit-level breakdown block (rows 17-19) is added, describing
#   how the Balancing-state bytes pack two 12-bit "Bal Cell" + "Temp"/"?"
#   fields per the 1A555408 message.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 changes
# ---------------------------------------------------------------------

# New label next to the CMC8 temperature row.
$ws1.Range("E65").Value = "And Balancing state"

# Rows 66-75: the "Balancing Request n CMC n" text currently in column B
# moves to column H, column B is cleared.
$rowsShiftToH = 66..75
foreach ($r in $rowsShiftToH) {
    $cellB = $ws1.Cells.Item($r, 2)
    $val = $cellB.Value()
    $ws1.Cells.Item($r, 8).Value = $val
    $cellB.ClearContents()
}

# Row 76 & 77: old column-B value moves to column H, column B is
# rewritten with the corrected (CMC 7) request text.
$ws1.Range("H76").Value = $ws1.Range("B76").Value()
$ws1.Range("B76").Value = "Balancing Request 1 CMC 7"

$ws1.Range("H77").Value = $ws1.Range("B77").Value()
$ws1.Range("B77").Value = "Balancing Request 2 CMC 7"

# Row 78 & 79: column B corrected to CMC 8 text (no column H value).
$ws1.Range("B78").Value = "Balancing Request 1 CMC 8"
$ws1.Range("B79").Value = "Balancing Request 2 CMC 8"

# Rows 82 & 83: the old "Balancing Request n CMC 8" text in column B is
# removed outright (it now lives nowhere, duplicated data cleaned up).
$ws1.Range("B82").ClearContents()
$ws1.Range("B83").ClearContents()

# ---------------------------------------------------------------------
# Sheet2 changes - add the new bit-level breakdown (rows 17-19)
# ---------------------------------------------------------------------

# Row 17: header byte indices, reusing the same layout/style as row 3.
$ws2.Range("A3").Copy()
$ws2.Range("A17").PasteSpecial(-4122)
$ws2.Range("A17").Value = "ID"

$ws2.Range("B3").Copy()
$headerCols = @("B", "C", "D", "L", "T", "U", "V", "W")
foreach ($col in $headerCols) {
    $ws2.Range($col + "17").PasteSpecial(-4122)
}
$ws2.Range("B17").Value = 0
$ws2.Range("C17").Value = 1
$ws2.Range("D17").Value = 2
$ws2.Range("L17").Value = 3
$ws2.Range("T17").Value = 4
$ws2.Range("U17").Value = 5
$ws2.Range("V17").Value = 6
$ws2.Range("W17").Value = 7

# Row 18: per-bit numbering (7..0 repeated across the two 8-bit fields).
$ws2.Range("B2").Copy()
$bitCols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S")
foreach ($col in $bitCols) {
    $ws2.Range($col + "18").PasteSpecial(-4122)
}
$bitValues = @(7, 6, 5, 4, 3, 2, 1, 0, 7, 6, 5, 4, 3, 2, 1, 0)
for ($i = 0; $i -lt $bitCols.Length; $i++) {
    $ws2.Range($bitCols[$i] + "18").Value = $bitValues[$i]
}

# Row 19: the actual field map for message 1A555408 ("Temp" / "?" bytes
# followed by the bit-packed Bal Cell values).
$ws1.Range("A65").Copy()
$ws2.Range("A19").PasteSpecial(-4122)
$ws2.Range("A19").Value = "1A555408"

$ws2.Range("B2").Copy()
$ws2.Range("B19").PasteSpecial(-4122)
$ws2.Range("C19").PasteSpecial(-4122)
$ws2.Range("B19").Value = "Temp"
$ws2.Range("C19").Value = "?"

$unstyledCells19 = @("D19", "E19", "F19", "G19", "H19", "I19", "J19", "K19", "P19", "Q19", "R19", "S19")
$unstyledValues19 = @("Bal Cell 8", "Bal Cell 7", "Bal Cell 6", "Bal Cell 5", "Bal Cell 4", "Bal Cell 3", "Bal Cell 2", "Bal Cell 1", "Bal Cell 12", "Bal Cell 11", "Bal Cell 10", "Bal Cell 9")
for ($i = 0; $i -lt $unstyledCells19.Length; $i++) {
    $cell = $ws2.Range($unstyledCells19[$i])
    $cell.Value = $unstyledValues19[$i]
    # These cells carry no explicit style in the source file - force off the
    # auto-inherited neighbour style so they serialise unstyled.
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# View state - match the author's final selection/scroll position.
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 57
$ws1.Range("D57").Select()

$ws2.Activate()
$ws2.Range("P22").Select()
